$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/11/2023  Through  9/17/2023"

# --- Simple numeric value updates (style unchanged) ---
$ws.Range("N14").Value = -88.888888888888
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -66.666666666666
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = 20
$ws.Range("L15").Value = -25
$ws.Range("N15").Value = -79.661016949152
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 96
$ws.Range("J16").Value = 123
$ws.Range("K16").Value = -21.951219512195
$ws.Range("L16").Value = 18.518518518518
$ws.Range("M16").Value = -50.259067357513
$ws.Range("N16").Value = -92.026578073089
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -22.222222222222
$ws.Range("F17").Value = 40
$ws.Range("G17").Value = 33
$ws.Range("H17").Value = 21.212121212121
$ws.Range("I17").Value = 255
$ws.Range("J17").Value = 249
$ws.Range("K17").Value = 2.409638554216
$ws.Range("L17").Value = 25
$ws.Range("M17").Value = 18.055555555555
$ws.Range("N17").Value = -63.150289017341
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 16
$ws.Range("H18").Value = -37.5
$ws.Range("I18").Value = 55
$ws.Range("J18").Value = 99
$ws.Range("K18").Value = -44.444444444444
$ws.Range("L18").Value = -6.779661016949
$ws.Range("M18").Value = -75.446428571428
$ws.Range("N18").Value = -95.898583146905
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 250
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = 100
$ws.Range("I19").Value = 292
$ws.Range("J19").Value = 244
$ws.Range("K19").Value = 19.672131147541
$ws.Range("L19").Value = 59.562841530054
$ws.Range("M19").Value = -1.683501683501
$ws.Range("N19").Value = -40.041067761807
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = -100
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = -26
$ws.Range("M20").Value = -35.087719298245
$ws.Range("N20").Value = -92.315680166147
$ws.Range("C21").Value = 18
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 117
$ws.Range("G21").Value = 93
$ws.Range("H21").Value = 25.806451612903
$ws.Range("I21").Value = 787
$ws.Range("J21").Value = 830
$ws.Range("K21").Value = -5.180722891566
$ws.Range("L21").Value = 24.920634920634
$ws.Range("M21").Value = -25.614366729678
$ws.Range("N21").Value = -83.511418395139
$ws.Range("M22").Value = -62.857142857142
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 4.761904761904
$ws.Range("F24").Value = 90
$ws.Range("H24").Value = 8.433734939759
$ws.Range("I24").Value = 801
$ws.Range("J24").Value = 793
$ws.Range("K24").Value = 1.008827238335
$ws.Range("L24").Value = 7.806191117092
$ws.Range("M24").Value = 11.871508379888
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 48
$ws.Range("G25").Value = 48
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 370
$ws.Range("J25").Value = 426
$ws.Range("K25").Value = -13.145539906103
$ws.Range("L25").Value = 23.333333333333
$ws.Range("M25").Value = -36.535162950257
$ws.Range("F26").Value = 1
$ws.Range("H26").Value = -66.666666666666
$ws.Range("J26").Value = 18
$ws.Range("K26").Value = -16.666666666666
$ws.Range("L26").Value = -34.782608695652
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = 100
$ws.Range("I27").Value = 41
$ws.Range("K27").Value = 2.5
$ws.Range("L27").Value = -2.380952380952
$ws.Range("F28").Value = 1
$ws.Range("L28").Value = -45
$ws.Range("F29").Value = 1
$ws.Range("L29").Value = -41.176470588235
$ws.Range("I30").Value = 9
$ws.Range("K30").Value = 12.5
$ws.Range("L30").Value = 200

# --- Numeric updates that also require a style/number-format change ---
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C18").Value = 2
$ws.Range("C18").NumberFormat = '#,##0'
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = '#,##0'
$ws.Range("E26").Value = -100
$ws.Range("E26").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = '#,##0'
$ws.Range("C30").Value = 1
$ws.Range("C30").NumberFormat = '#,##0'
$ws.Range("F30").Value = 1
$ws.Range("F30").NumberFormat = '#,##0'

# --- Cells that become the "no data" text placeholder (copy format+value from a matching cell) ---
$ws.Range("D28").Copy($ws.Range("C20"))
$ws.Range("D28").Copy($ws.Range("C28"))
$ws.Range("D28").Copy($ws.Range("C29"))
